# Update Tnfsf18-Tnfrsf18.xlsx with refreshed TPM-derived values.
# Rows 5-7 (the old "FAPs -> *" block) are removed, and rows 2-4
# (the old "ECs -> *" block) are rewritten in-place with the numbers
# that used to live in rows 5-7, recomputed against the new TPM totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-redundant trailing rows first.
$ws.Range("A5:T7").EntireRow.Delete()

function Set-RowValues {
    param($ws, $row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Row 2: FAPs -> ECs
Set-RowValues $ws 2 @("FAPs", "Tnfsf18", "Tnfrsf18", "ECs", 2, 0.6666666666666666, 0.3800126666666667, 1.140038, 1, 1, 3, 1, 0.4588346666666667, 1.376504, 0.1133721507248096, 0.1133721507248096, 0.1743629852391111, 1.569266867152, 0.1133721507248096, 0.1133721507248096)

# Row 3: FAPs -> FAPs
Set-RowValues $ws 3 @("FAPs", "Tnfsf18", "Tnfrsf18", "FAPs", 2, 0.6666666666666666, 0.3800126666666667, 1.140038, 1, 1, 3, 1, 2.108323666666667, 6.324971, 0.5209396889090402, 0.5209396889090402, 0.8011896987664445, 7.210707288898, 0.5209396889090402, 0.5209396889090402)

# Row 4: FAPs -> MuSCs
Set-RowValues $ws 4 @("FAPs", "Tnfsf18", "Tnfrsf18", "MuSCs", 2, 0.6666666666666666, 0.3800126666666667, 1.140038, 1, 1, 3, 1, 1.479996666666667, 4.43999, 0.3656881603661502, 0.3656881603661502, 0.5624174799577778, 5.061757319620001, 0.3656881603661502, 0.3656881603661502)
